$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18..107 down to 19..108)
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new price-report record.
# Columns A, B, C, E, F, G, H, I, J, R keep the same "template" values used
# throughout this sheet for Vega Monumental Concepción / Bíobío / Durazno.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44537
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100103
$ws.Range("H18").Value = "Frutos de hueso (carozo)"
$ws.Range("I18").Value = 100103004
$ws.Range("J18").Value = "Durazno"
$ws.Range("K18").Value = "Early Majestic"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 250
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15520
$ws.Range("Q18").Value = "$/caja 15 kilos empedrada"
$ws.Range("R18").Value = "Región de O'Higgins"
$ws.Range("S18").Value = 1035
$ws.Range("T18").Value = 15
